$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# --- Slide 1: Subtitle "By Ethan Laviolette" -> runs "Ethan " + "Laviolette" ---
$s1 = $p.Slides.Item(1)
$subtitle = Get-ShapeByName $s1 "Subtitle 2"
$tr1 = $subtitle.TextFrame.TextRange
$tr1.Text = "Ethan Laviolette"
$front = $subtitle.TextFrame.TextRange.Characters(1, 6)
$front.Font.Bold = $true

# --- Slide 5 (last slide): Title "	The " + "Results" -> single run "	The Results" ---
$s5 = $p.Slides.Item($p.Slides.Count)
$title = Get-ShapeByName $s5 "Title 1"
$tr5 = $title.TextFrame.TextRange
$tr5.Characters(1, 5).Delete()
$remaining = $title.TextFrame.TextRange
$remaining.Text = "`tThe Results"
